$wb = $excel.ActiveWorkbook

# --- Sheet "rpa" (sheet1): update BTC price and selection ---
$ws1 = $wb.Worksheets.Item("rpa")
$ws1.Range("B2").Value = 40000

# --- Sheet "buyValue" (sheet2): append new price/percent log rows ---
$ws2 = $wb.Worksheets.Item("buyValue")

$rows = @(
    @("XRP", 0.3125),
    @("XRP", 0.3179),
    @("XRP", 0.3138),
    @("BTC", 38662.12),
    @("XRP", 0.311),
    @("BTC", 34713.25),
    @("XRP", 0.292),
    @("BTC", 34609.23),
    @("XRP", 0.2922),
    @("BTC", 34257.21),
    @("BTC", 34286.88),
    @("BTC", 34029.62),
    @("XRP", 0.2906),
    @("BTC", 34036.73),
    @("XRP", 0.2913)
)

$r = 9
foreach ($row in $rows) {
    $ws2.Range("A$r").Value = $row[0]
    $ws2.Range("B$r").Value = $row[1]
    $r = $r + 1
}

# --- Selections: set sheet1's selection first, then sheet2's last so
#     "buyValue" remains the active tab (as in the original workbook) ---
$ws1.Range("D5").Select() | Out-Null
$ws2.Range("B16").Select() | Out-Null
